# Fixing automation script for suite E
# Applies the B-suite row additions/edits to the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# Column layout: column B gets its own (wider) width, split off from column A
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 19.5

# ---------------------------------------------------------------------------
# Row 100: D keeps "Y" but loses the (redundant) fill flag; E becomes SKIP
# ---------------------------------------------------------------------------
$ws.Range("A101").Copy()
$ws.Range("D100").PasteSpecial($xlPasteFormats)
$ws.Range("D100").Value = "Y"
$ws.Range("E100").Value = "SKIP"

# ---------------------------------------------------------------------------
# Rows 101-104: E column goes from blank to SKIP (style unchanged, s=3)
# ---------------------------------------------------------------------------
$ws.Range("E101").Value = "SKIP"
$ws.Range("E102").Value = "SKIP"
$ws.Range("E103").Value = "SKIP"
$ws.Range("E104").Value = "SKIP"

# ---------------------------------------------------------------------------
# Row 105: C becomes the "ALL search results" verification text, D keeps Y
# (losing the redundant fill flag), E becomes SKIP
# ---------------------------------------------------------------------------
$ws.Range("C105").Value = "Verify that record view page of a post gets displayed when user clicks on article title in ALL  search results page"
$ws.Range("A101").Copy()
$ws.Range("D105").PasteSpecial($xlPasteFormats)
$ws.Range("D105").Value = "Y"
$ws.Range("E105").Value = "SKIP"

# ---------------------------------------------------------------------------
# Row 106: B becomes the combined OPQA ids, C becomes the long combined
# multi-line verification text (now wrap-formatted + taller row), E -> SKIP
# ---------------------------------------------------------------------------
$ws.Range("B106").Value = "OPQA-555|OPQA-556"

$c106 = @(
  "Verify that record view page of a post gets displayed when user clicks on article title in POSTs search results page",
  "Verify that following fields get displayed correctly for a post in record view page:",
  "a)Title",
  " b)Creation date and time ",
  "c)Last edited date and time ",
  "d)Author ",
  "e)Author details ",
  "f)Post content ",
  "g)Likes count ",
  "h)Comments count ",
  "i)Views count"
) -join "`r`n"

$ws.Range("C27").Copy()
$ws.Range("C106").PasteSpecial($xlPasteFormats)
$ws.Range("C106").Value = $c106

$ws.Range("E106").Value = "SKIP"
$ws.Rows.Item(106).RowHeight = 165

# ---------------------------------------------------------------------------
# New row 107 (TestCase_B106 / OPQA-1226 / SORT BY options for POSTS)
# ---------------------------------------------------------------------------
$ws.Range("A101:E101").Copy()
$ws.Range("A107:E107").PasteSpecial($xlPasteFormats)

$ws.Range("A107").Value = "TestCase_B106"
$ws.Range("B107").Value = "OPQA-1226"

$c107 = @(
  "Verify that following options get displayed in SORT BY drop down in POSTS search results page: ",
  "a)Relevance ",
  "b)Create Date(Newest) ",
  "c)Create Date(Oldest)"
) -join "`r`n"
$ws.Range("C27").Copy()
$ws.Range("C107").PasteSpecial($xlPasteFormats)
$ws.Range("C107").Value = $c107

$ws.Range("D107").Value = "Y"
$ws.Range("E107").Value = "PASS"

# B107 gets the new black-font style (fontId=2 in the finished workbook)
$ws.Range("B107").Font.Color = 0

$ws.Rows.Item(107).RowHeight = 60

# ---------------------------------------------------------------------------
# New row 108 (TestCase_B107 / OPQA-574 / left nav retained - ALL)
# ---------------------------------------------------------------------------
$ws.Range("A71:E71").Copy()
$ws.Range("A108:E108").PasteSpecial($xlPasteFormats)

$ws.Range("A108").Value = "TestCase_B107"
$ws.Range("B108").Value = "OPQA-574"
$ws.Range("C108").Value = "Verify that left navigation pane content type is retained when user navigates back to ALL search results page from record view page"
$ws.Range("D108").Value = "Y"
$ws.Range("E108").Value = "SKIP"

# ---------------------------------------------------------------------------
# New row 109 (TestCase_B108 / OPQA-569 / sorting retained - ALL)
# ---------------------------------------------------------------------------
$ws.Range("A71:E71").Copy()
$ws.Range("A109:E109").PasteSpecial($xlPasteFormats)

$ws.Range("A109").Value = "TestCase_B108"
$ws.Range("B109").Value = "OPQA-569"
$ws.Range("C109").Value = "Verify that sorting is retained when user navigates back to ALL search results page from record view page"
$ws.Range("D109").Value = "Y"
$ws.Range("E109").Value = "PASS"

# ---------------------------------------------------------------------------
# New row 110 (TestCase_B109 / OPQA-592 / PATENTS filters)
# ---------------------------------------------------------------------------
$ws.Range("A101:E101").Copy()
$ws.Range("A110:E110").PasteSpecial($xlPasteFormats)

$ws.Range("A110").Value = "TestCase_B109"
$ws.Range("B110").Value = "OPQA-592"
$ws.Range("C110").Value = "Verify that following filters are present in PATENTS search results page: a)Inventor b)IPC Codes c)Assignee"
$ws.Range("D110").Value = "Y"

# ---------------------------------------------------------------------------
# New row 111 (TestCase_B110 / OPQA-581 / TIMES CITED sort - PATENTS)
# ---------------------------------------------------------------------------
$ws.Range("A101:E101").Copy()
$ws.Range("A111:E111").PasteSpecial($xlPasteFormats)

$ws.Range("A111").Value = "TestCase_B110"
$ws.Range("B111").Value = "OPQA-581"
$ws.Range("C111").Value = "Verify that search results are sorted correctly by TIMES CITED field in SORT BY drop down in PATENTS search results page"
$ws.Range("D111").Value = "Y"

# ---------------------------------------------------------------------------
# New row 112 (TestCase_B111 / OPQA-1242 / PEOPLE scroll more results)
# ---------------------------------------------------------------------------
$ws.Range("A100:E100").Copy()
$ws.Range("A112:E112").PasteSpecial($xlPasteFormats)

$ws.Range("A112").Value = "TestCase_B111"
$ws.Range("B112").Value = "OPQA-1242"
$ws.Range("C112").Value = "Verify that more search results get displayed when user scrolls down in PEOPLE search results page."
$ws.Range("D112").Value = "Y"
$ws.Range("E112").Value = "SKIP"

# ---------------------------------------------------------------------------
# New row 113 (TestCase_B112 / OPQA-1243 / sorting retained - PEOPLE)
# ---------------------------------------------------------------------------
$ws.Range("A100:E100").Copy()
$ws.Range("A113:E113").PasteSpecial($xlPasteFormats)

$ws.Range("A113").Value = "TestCase_B112"
$ws.Range("B113").Value = "OPQA-1243"
$ws.Range("C113").Value = "Verify that sorting is retained when user navigates back to PEOPLE search results page from record view page"
$ws.Range("D113").Value = "Y"
$ws.Range("E113").Value = "SKIP"

# ---------------------------------------------------------------------------
# New rows 114-115: plain/unformatted rows (no explicit cell style), matching
# TestCase_B113 / OPQA-593 / INVENTOR MORE-LESS and TestCase_B114 / OPQA-588
# ---------------------------------------------------------------------------
$ws.Range("A114").Value = "TestCase_B113"
$ws.Range("B114").Value = "OPQA-593"
$ws.Range("C114").Value = "Verify that MORE and LESS links are working correctly in INVENTOR filter in PATENTS search results page"
$ws.Range("D114").Value = "Y"
$ws.Range("E114").Value = "PASS"

$ws.Range("A115").Value = "TestCase_B114"
$ws.Range("B115").Value = "OPQA-588"
$ws.Range("C115").Value = "Verify that left navigation pane content type is retained when user navigates back to PATENTS search results page from record view page"
$ws.Range("D115").Value = "Y"
$ws.Range("E115").Value = "PASS"

# ---------------------------------------------------------------------------
# View state: scroll to the bottom, select the last new row
# ---------------------------------------------------------------------------
$ws.Range("A115").Select()
$excel.ActiveWindow.ScrollRow = 107

Write-Output "edit complete"
